$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- India: updated daily figures (row 19) ---
$ws.Range("B19").Value = 39980
$ws.Range("C19").Value = 281
$ws.Range("E19").Value = 27838

# --- Reorder Kirguistan / Albania (rows 96-97) and refresh their daily figures ---
# Kirguistan now sits above Albania, each carrying its own updated counts.
$ws.Range("A96").Value = "Kirguistan"
$ws.Range("B96").Value = 795
$ws.Range("C96").Value = 26
$ws.Range("D96").Value = 564
$ws.Range("E96").Value = 221
$ws.Range("F96").Value = 12
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 10

$ws.Range("A97").Value = "Albania"
$ws.Range("B97").Value = 789
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 519
$ws.Range("E97").Value = 239
$ws.Range("F97").Value = 4
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 31

# --- Reorder Seychelles / Montserrat (rows 204-205) and refresh their daily figures ---
$ws.Range("A204").Value = "Seychelles"
$ws.Range("B204").Value = 11
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 6
$ws.Range("E204").Value = 5
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

$ws.Range("A205").Value = "Montserrat"
$ws.Range("B205").Value = 11
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 7
$ws.Range("E205").Value = 3
$ws.Range("F205").Value = 1
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 1
